# Update "想去人数" (interest count) values in the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 6830
$ws1.Range("F4").Value = 437
$ws1.Range("F8").Value = 110
$ws1.Range("F9").Value = 104
$ws1.Range("F12").Value = 30
$ws1.Range("F14").Value = 425
$ws1.Range("F15").Value = 4
$ws1.Range("F18").Value = 3447
$ws1.Range("F21").Value = 12
$ws1.Range("F22").Value = 2082
$ws1.Range("F23").Value = 172
$ws1.Range("F26").Value = 3

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6830
$ws4.Range("F4").Value = 437
$ws4.Range("F9").Value = 110
$ws4.Range("F10").Value = 104
$ws4.Range("F13").Value = 30
$ws4.Range("F15").Value = 425
$ws4.Range("F16").Value = 4
$ws4.Range("F19").Value = 3447
$ws4.Range("F22").Value = 12
$ws4.Range("F23").Value = 2082
$ws4.Range("F24").Value = 172
$ws4.Range("F27").Value = 3
